# Report types.xlsx — add the "coreMetricsBySite" / "core-metrics-table-section"
# row to the Real User report's test-case list (new row 21, columns A:B),
# matching the existing formatting used by the row above it (A20:B20).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Carry the existing A20:B20 formatting down to A21:B21 (same style indices
# used throughout column A / column B of this table) before filling values.
$ws.Range("A20:B20").Copy()
$ws.Range("A21:B21").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A21").Value = "coreMetricsBySite"
$ws.Range("B21").Value = "core-metrics-table-section"

# Leave the selection on B21, as in the saved workbook.
[void]$ws.Range("B21").Select()
